$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.391.90'
$ws.Range("E2").Value = '  -0.69%  '
$ws.Range("D3").Value = '1.639.11'
$ws.Range("E3").Value = '  -1.60%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.86'
$ws.Range("E5").Value = '  -2.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.533'
$ws.Range("E6").Value = '  +3.80%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.11'
$ws.Range("E8").Value = '  -2.24%  '
$ws.Range("E9").Value = '  -3.07%  '
$ws.Range("E10").Value = '  -2.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0890'
$ws.Range("E11").Value = '  +1.03%  '
$ws.Range("D12").Value = '1.871.02'
$ws.Range("E12").Value = '  -1.58%  '
$ws.Range("D13").Value = '1.640.57'
$ws.Range("E13").Value = '  -1.50%  '
$ws.Range("E14").Value = '  -2.73%  '
$ws.Range("E15").Value = '  -1.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.21'
$ws.Range("E16").Value = '  -2.97%  '
$ws.Range("D17").Value = '27.359.09'
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '229.43'
$ws.Range("E18").Value = '  -5.29%  '
$ws.Range("E19").Value = '  -1.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.50'
$ws.Range("E20").Value = '  -1.04%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("E22").Value = '  -4.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.33'
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("E24").Value = '  -1.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.07'
$ws.Range("E25").Value = '  +1.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.95'
$ws.Range("E26").Value = '  -3.39%  '
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("E29").Value = '  -5.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.17'
$ws.Range("E30").Value = '  -5.25%  '
$ws.Range("E31").Value = '  -3.96%  '
$ws.Range("E32").Value = '  -2.35%  '
$ws.Range("E33").Value = '  -0.16%  '
$ws.Range("D34").Value = '1.408.21'
$ws.Range("E34").Value = '  -4.47%  '
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("E37").Value = '  -2.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.879'
$ws.Range("E38").Value = '  -5.80%  '
$ws.Range("E39").Value = '  -2.99%  '
$ws.Range("E40").Value = '  +0.90%  '
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.47'
$ws.Range("E42").Value = '  -1.83%  '
$ws.Range("E43").Value = '  +0.75%  '
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("E45").Value = '  +0.20%  '
$ws.Range("E46").Value = '  -7.23%  '
$ws.Range("D47").Value = '1.780.93'
$ws.Range("E47").Value = '  -1.51%  '
$ws.Range("E48").Value = '  -4.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '87.24'
$ws.Range("E49").Value = '  -2.44%  '
$ws.Range("D50").Value = '0.0₆0105'
$ws.Range("E50").Value = '  -2.34%  '
$ws.Range("E51").Value = '  -3.83%  '
